# Shopping slots export — add agency + status columns, rename/reorder
# fields, refresh the sample rows into the real exported slot rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column layout changes
#    - insert "대행사" (agency) as the new column B
#    - insert "상태" (status) as a new column right before the last
#      column ("슬롯 단가" / slot unit price)
# ---------------------------------------------------------------------
$ws.Columns("B").Insert()
$ws.Columns("L").Insert()

# ---------------------------------------------------------------------
# 2. Header row (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "슬롯명"
$ws.Range("B1").Value = "대행사"
$ws.Range("C1").Value = "스토어 타입"
$ws.Range("D1").Value = "상품 ID"
$ws.Range("E1").Value = "상품명"
$ws.Range("F1").Value = "키워드"
$ws.Range("G1").Value = "가격"
$ws.Range("H1").Value = "할인가"
$ws.Range("I1").Value = "시작일"
$ws.Range("J1").Value = "종료일"
$ws.Range("K1").Value = "입찰방식"
$ws.Range("L1").Value = "상태"
$ws.Range("M1").Value = "슬롯 단가"

# ---------------------------------------------------------------------
# 3. Replace the two sample rows with the ten live export rows.
#    Clear out the old sample data first (columns C:K of the old
#    layout, now C:M minus the two new columns) ...
# ---------------------------------------------------------------------
$ws.Range("A2:M3").ClearContents()

# ... then add eight more rows so we have ten data rows (2-11).
$ws.Rows("4:11").Insert()

$startDate = 45784
$endDate = 45786

for ($i = 1; $i -le 10; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = "이이이 Shopping 슬롯 20250506-$i"
    $ws.Cells.Item($r, 2).Value = "이이이"
    $ws.Cells.Item($r, 9).Value = $startDate
    $ws.Cells.Item($r, 9).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 10).Value = $endDate
    $ws.Cells.Item($r, 10).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($r, 12).Value = "live"
    $ws.Cells.Item($r, 13).Value = 30
}

Write-Output "done"
